$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.459.44"
$ws.Range("E2").Value = "  -1.11%  "

$ws.Range("D3").Value = "2.626.40"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("E4").Value = "  -0.04%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.21"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -2.36%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.01"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.00%  "

$ws.Range("E7").Value = "  -1.58%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -2.38%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.34"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -4.74%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.09"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -0.42%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0808"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("E14").Value = "  -0.43%  "

$ws.Range("D15").Value = "3.036.80"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("D16").Value = "2.641.42"
$ws.Range("E16").Value = "  -0.90%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -3.12%  "

$ws.Range("D18").Value = "49.386.45"
$ws.Range("E18").Value = "  -1.18%  "

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.66"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.89"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.56%  "

$ws.Range("E22").Value = "  -1.71%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.96"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -4.04%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.79"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -5.06%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("E26").Value = "  +0.14%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.98"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -3.54%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +1.55%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.15"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -3.11%  "

$ws.Range("E30").Value = "  -2.04%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.52"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -5.06%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.59"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -1.71%  "

$ws.Range("E33").Value = "  +0.58%  "

$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("E36").Value = "  -3.33%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +3.12%  "

$ws.Range("E38").Value = "  -3.28%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -0.17%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "129.02"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +2.40%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.85"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -0.87%  "

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0324"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +2.92%  "

$ws.Range("D45").Value = "2.038.96"
$ws.Range("E45").Value = "  -1.81%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.17"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +8.50%  "

$ws.Range("E47").Value = "  -4.63%  "

$ws.Range("E48").Value = "  -3.94%  "

$ws.Range("E49").Value = "  -3.77%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.20"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -4.00%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.39"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.84%  "
